# Pins.xlsx — "Worked on adding pin numbers"
#
# Adds pin-out tables for a BNO055 IMU, a push Button, an SD card
# adapter, and an NRF24L01 module to the existing Microcontroller
# pinout legend, and relocates/restyles the "Key" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter  = -4108
$xlGeneral = 1

# ---------------------------------------------------------------
# 1. Break apart the merged ranges that are being reshaped so the
#    cells underneath become individually addressable again.
# ---------------------------------------------------------------
$ws.Range("A2:A20").UnMerge()
$ws.Range("C2:E2").UnMerge()
$ws.Range("C35:D35").UnMerge()
$ws.Range("C36:D39").UnMerge()

# ---------------------------------------------------------------
# 2. Row 2 / Row 3 — accessory group headers and their pin labels.
# ---------------------------------------------------------------

# Row 2 group headers
$ws.Range("C2").Value = "BNO055IMU"
$ws.Range("I2").Value = "Button"
$ws.Range("K2").Value = "SD Adapter"
$ws.Range("Q2").Value = "NRF24l01 standard"

# Row 3 pin labels under each group
$ws.Range("C3").Value = "Vin"
$ws.Range("D3").Value = "3vo"
$ws.Range("E3").Value = "GND"
$ws.Range("F3").Value = "SDA"
$ws.Range("G3").Value = "SCL"
$ws.Range("H3").Value = "RST"
$ws.Range("I3").Value = "side1"
$ws.Range("J3").Value = "side2"
$ws.Range("K3").Value = "CS"
$ws.Range("L3").Value = "SCK"
$ws.Range("M3").Value = "MOSI"
$ws.Range("N3").Value = "MISO"
$ws.Range("O3").Value = "VCC"
$ws.Range("P3").Value = "GND"

# Apply the centred style used throughout the header band to the
# newly populated / widened header cells.
$ws.Range("C2:R3").HorizontalAlignment = $xlCenter
$ws.Range("C2:R3").VerticalAlignment = $xlCenter

# ---------------------------------------------------------------
# 3. Re-merge the header groups at their new extents.
# ---------------------------------------------------------------
$ws.Range("C2:H2").Merge()
$ws.Range("I2:J2").Merge()
$ws.Range("K2:P2").Merge()
$ws.Range("Q2:W2").Merge()

# "Microcontroller" label now only spans rows 2-6 and reads
# vertically (rotated 90 degrees).
$ws.Range("A2:A6").Merge()
$ws.Range("A2:A6").HorizontalAlignment = $xlCenter
$ws.Range("A2:A6").VerticalAlignment = $xlCenter
$ws.Range("A2:A6").Orientation = 90
$ws.Range("A2:A6").WrapText = $false

# ---------------------------------------------------------------
# 4. Rows that used to fall inside A2:A20 but no longer sit under
#    the "Microcontroller" merge — drop back to a plain vertical
#    centering (no horizontal centering).
# ---------------------------------------------------------------
$ws.Range("A7:A20").HorizontalAlignment = $xlGeneral
$ws.Range("A7:A20").VerticalAlignment = $xlCenter

# ---------------------------------------------------------------
# 5. "Key" legend block — remove the old label from C35 and place
#    it, rotated, in the taller merged C37:D39 cell instead.
# ---------------------------------------------------------------
$ws.Range("C35").Value = ""
$ws.Range("C35:D35").HorizontalAlignment = $xlGeneral
$ws.Range("C35:D35").VerticalAlignment = $xlCenter

$ws.Range("C36:D36").Merge()
$ws.Range("C36:D36").HorizontalAlignment = $xlCenter
$ws.Range("C36:D36").VerticalAlignment = $xlCenter
$ws.Range("C36:D36").WrapText = $false

$ws.Range("C37:D39").Merge()
$ws.Range("C37").Value = "Key"
$ws.Range("C37:D39").HorizontalAlignment = $xlCenter
$ws.Range("C37:D39").VerticalAlignment = $xlCenter
$ws.Range("C37:D39").Orientation = 90
$ws.Range("C37:D39").WrapText = $true

# Rows 37-39 grow to fit the taller legend rows (pictures/links
# area gets more breathing room).
$ws.Rows(37).RowHeight = 25.5
$ws.Rows(38).RowHeight = 25.5
$ws.Rows(39).RowHeight = 25.5

# ---------------------------------------------------------------
# 6. Move the selection like the author left it.
# ---------------------------------------------------------------
$ws.Range("S3").Select()
